$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update data values (rows 36-37)
$ws.Range("C36").Value = 0.75
$ws.Range("A37").Value = 43.75
$ws.Range("C37").Value = 0.3

# Update the view: scroll so row 17 is the top-left row, and select E34
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 17
$ws.Range("E34").Select()
